# This workbook contains a weekly price-report dataset for "Papa" (potato)
# at "Femacal de La Calera". A new week of data (3 rows) needs to be
# inserted just above the existing row 283, pushing the old rows 283-355
# down to 286-358 (and updating the sheet dimension accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 283:285, shifting former rows 283-355 down to 286-358.
$ws.Rows("283:285").Insert()

# ---- Row 283 (new) ----
$ws.Cells.Item(283, 1).Value = 3
$ws.Cells.Item(283, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(283, 3).Value = "Coquimbo"
$ws.Cells.Item(283, 4).Value = 44508
$ws.Cells.Item(283, 5).Value = 5
$ws.Cells.Item(283, 6).Value = 100114001
$ws.Cells.Item(283, 7).Value = "Papa"
$ws.Cells.Item(283, 8).Value = "Asterix"
$ws.Cells.Item(283, 9).Value = "1a nueva(o)"
$ws.Cells.Item(283, 10).Value = 510
$ws.Cells.Item(283, 11).Value = 9000
$ws.Cells.Item(283, 12).Value = 9500
$ws.Cells.Item(283, 13).Value = 9245
$ws.Cells.Item(283, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(283, 15).Value = "Provincia de Talca"
$ws.Cells.Item(283, 16).Value = 370
$ws.Cells.Item(283, 17).Value = 25
$ws.Cells.Item(283, 18).Value = "Hortaliza"

# ---- Row 284 (new) ----
$ws.Cells.Item(284, 1).Value = 3
$ws.Cells.Item(284, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(284, 3).Value = "Coquimbo"
$ws.Cells.Item(284, 4).Value = 44508
$ws.Cells.Item(284, 5).Value = 5
$ws.Cells.Item(284, 6).Value = 100114001
$ws.Cells.Item(284, 7).Value = "Papa"
$ws.Cells.Item(284, 8).Value = "Rosara"
$ws.Cells.Item(284, 9).Value = "1a (cosecha)"
$ws.Cells.Item(284, 10).Value = 260
$ws.Cells.Item(284, 11).Value = 9000
$ws.Cells.Item(284, 12).Value = 9000
$ws.Cells.Item(284, 13).Value = 9000
$ws.Cells.Item(284, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(284, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(284, 16).Value = 360
$ws.Cells.Item(284, 17).Value = 25
$ws.Cells.Item(284, 18).Value = "Hortaliza"

# ---- Row 285 (new) ----
$ws.Cells.Item(285, 1).Value = 3
$ws.Cells.Item(285, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(285, 3).Value = "Coquimbo"
$ws.Cells.Item(285, 4).Value = 44508
$ws.Cells.Item(285, 5).Value = 5
$ws.Cells.Item(285, 6).Value = 100114001
$ws.Cells.Item(285, 7).Value = "Papa"
$ws.Cells.Item(285, 8).Value = "Rosara"
$ws.Cells.Item(285, 9).Value = "1a nueva(o)"
$ws.Cells.Item(285, 10).Value = 250
$ws.Cells.Item(285, 11).Value = 9500
$ws.Cells.Item(285, 12).Value = 9500
$ws.Cells.Item(285, 13).Value = 9500
$ws.Cells.Item(285, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(285, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(285, 16).Value = 380
$ws.Cells.Item(285, 17).Value = 25
$ws.Cells.Item(285, 18).Value = "Hortaliza"

Write-Host ("Dimension: " + $ws.Range("A1:R358").Address())
